# Update the "Förändrad" (changed) date column (C) for rows 2-23
# from serial date 45183 to 45184 (one day later), matching the
# automatic update performed upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
